$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '76.451.83'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.880.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.82%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '196.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '599.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.554'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.67%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.881.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.78%  '
$ws.Range("E11").Value = '  +9.84%  '
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("E13").Value = '  +4.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.411.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.368.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.36%  '
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.867.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.89%  '
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.029.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.88%  '
$ws.Range("E29").Value = '  +10.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '511.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("E34").Value = '  +3.46%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '167.17'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.04'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.46%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '185.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.19%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("E42").Value = '  +4.61%  '
$ws.Range("E43").Value = '  +1.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0919'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.49%  '
$ws.Range("E46").Value = '  +3.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.580'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.680'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +14.86%  '
$ws.Range("E51").Value = '  +3.30%  '
